$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update product_name column (B) for the "parka" bulk-upload rows
$ws.Range("B2").Value = "parka1"
$ws.Range("B3").Value = "parka2"
$ws.Range("B4").Value = "parka3"

# Update size column (M) from L to S
$ws.Range("M2").Value = "S"
$ws.Range("M3").Value = "S"
$ws.Range("M4").Value = "S"

# Update the view: scroll back to A1 (clear topLeftCell) and move the
# active selection to M4
$ws.Range("A1").Select()
$ws.Range("M4").Select()
